# The deck currently uses the "Integral" theme (green colour scheme) as its
# active theme. The edit re-colours the presentation to use the classic
# "Office Theme" colour scheme instead (dk2/lt2/accent1-6/hlink/folHlink),
# matching the target OOXML colour values exactly. dk1/lt1 are unchanged
# (already black/white in both schemes) but are (re)applied for completeness.
#
# PowerPoint COM RGB values are packed as 0x00BBGGRR (R + G*256 + B*65536),
# so each target hex colour below is converted accordingly.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

# index -> (name, target hex) taken from the "Office Theme" colour scheme
# 1 dk1, 2 lt1, 3 dk2, 4 lt2, 5 accent1, 6 accent2, 7 accent3, 8 accent4,
# 9 accent5, 10 accent6, 11 hlink, 12 folHlink
$tcs.Colors(1).RGB  = 0        # dk1      000000
$tcs.Colors(2).RGB  = 16777215 # lt1      FFFFFF
$tcs.Colors(3).RGB  = 6968388  # dk2      44546A
$tcs.Colors(4).RGB  = 15132391 # lt2      E7E6E6
$tcs.Colors(5).RGB  = 13998939 # accent1  5B9BD5
$tcs.Colors(6).RGB  = 3243501  # accent2  ED7D31
$tcs.Colors(7).RGB  = 10855845 # accent3  A5A5A5
$tcs.Colors(8).RGB  = 49407    # accent4  FFC000
$tcs.Colors(9).RGB  = 12874308 # accent5  4472C4
$tcs.Colors(10).RGB = 4697456  # accent6  70AD47
$tcs.Colors(11).RGB = 12673797 # hlink    0563C1
$tcs.Colors(12).RGB = 7491477  # folHlink 954F72
